# Dodato ciscenje starih logova
# Refresh the hydro-plant log sheet: update the "Time" column for the
# rows that remain, and drop the three oldest trailing log rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamps for the surviving log entries (rows 2-23, "No." 1-22).
$times = @(
    "08/01/2023 9:52:55 pm",
    "08/01/2023 9:52:56 pm",
    "08/01/2023 9:52:56 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:57 pm",
    "08/01/2023 9:52:58 pm",
    "08/01/2023 9:52:58 pm",
    "08/01/2023 9:52:58 pm",
    "08/01/2023 9:52:58 pm",
    "08/01/2023 9:52:58 pm",
    "08/01/2023 9:52:59 pm",
    "08/01/2023 9:52:59 pm",
    "08/01/2023 9:52:59 pm",
    "08/01/2023 9:53:01 pm",
    "08/01/2023 9:53:01 pm",
    "08/01/2023 9:53:02 pm",
    "08/01/2023 9:53:02 pm"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $times[$i]
}

# Remove the three oldest trailing log rows (old No. 23, 24, 25).
$ws.Rows("24:26").Delete()
